$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "56.355.22"
$ws.Cells.Item(2, 5).Value = "  -4.31%  "
$ws.Cells.Item(3, 4).Value = "2.368.83"
$ws.Cells.Item(3, 5).Value = "  -5.25%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "498.12"
$ws.Cells.Item(5, 5).Value = "  -6.56%  "
$ws.Cells.Item(6, 4).Value = "128.72"
$ws.Cells.Item(6, 5).Value = "  -3.80%  "
$ws.Cells.Item(7, 5).Value = "  -0.29%  "
$ws.Cells.Item(8, 4).Value = "'0.550"
$ws.Cells.Item(8, 5).Value = "  -3.39%  "
$ws.Cells.Item(9, 4).Value = "2.389.16"
$ws.Cells.Item(9, 5).Value = "  -4.51%  "
$ws.Cells.Item(10, 4).Value = "0.0954"
$ws.Cells.Item(10, 5).Value = "  -3.85%  "
$ws.Cells.Item(11, 5).Value = "  -1.39%  "
$ws.Cells.Item(12, 5).Value = "  -3.27%  "
$ws.Cells.Item(13, 4).Value = "4.59"
$ws.Cells.Item(13, 5).Value = "  -10.82%  "
$ws.Cells.Item(14, 4).Value = "2.793.53"
$ws.Cells.Item(14, 5).Value = "  -5.11%  "
$ws.Cells.Item(15, 4).Value = "56.847.32"
$ws.Cells.Item(15, 5).Value = "  -3.15%  "
$ws.Cells.Item(16, 4).Value = "21.42"
$ws.Cells.Item(16, 5).Value = "  -3.93%  "
$ws.Cells.Item(17, 5).Value = "  -3.76%  "
$ws.Cells.Item(18, 4).Value = "2.404.19"
$ws.Cells.Item(18, 5).Value = "  -3.83%  "
$ws.Cells.Item(19, 4).Value = "'10.10"
$ws.Cells.Item(19, 5).Value = "  -4.58%  "
$ws.Cells.Item(20, 4).Value = "'311.30"
$ws.Cells.Item(20, 5).Value = "  -2.99%  "
$ws.Cells.Item(21, 4).Value = "4.02"
$ws.Cells.Item(21, 5).Value = "  -5.40%  "
$ws.Cells.Item(22, 4).Value = "6.17"
$ws.Cells.Item(22, 5).Value = "  +0.07%  "
$ws.Cells.Item(23, 4).Value = "0.998"
$ws.Cells.Item(23, 5).Value = "  -0.09%  "
$ws.Cells.Item(24, 4).Value = "65.39"
$ws.Cells.Item(24, 5).Value = "  -0.58%  "
$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 5).Value = "  +0.18%  "
$ws.Cells.Item(26, 4).Value = "2.491.94"
$ws.Cells.Item(26, 5).Value = "  -4.96%  "
$ws.Cells.Item(27, 5).Value = "  -9.37%  "
$ws.Cells.Item(28, 4).Value = "0.149"
$ws.Cells.Item(28, 5).Value = "  -6.18%  "
$ws.Cells.Item(29, 4).Value = "7.17"
$ws.Cells.Item(29, 5).Value = "  -3.44%  "
$ws.Cells.Item(30, 4).Value = "173.99"
$ws.Cells.Item(30, 5).Value = "  +0.79%  "
$ws.Cells.Item(31, 5).Value = "  -4.66%  "
$ws.Cells.Item(32, 4).Value = "0.0₃0706"
$ws.Cells.Item(32, 5).Value = "  -6.35%  "
$ws.Cells.Item(33, 4).Value = "6.07"
$ws.Cells.Item(33, 5).Value = "  -3.12%  "
$ws.Cells.Item(34, 2).Value = "USDe"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(34, 4).Value = "0.998"
$ws.Cells.Item(34, 5).Value = "  -0.05%  "
$ws.Cells.Item(35, 2).Value = "Fetch.AI"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(35, 4).Value = "1.09"
$ws.Cells.Item(35, 5).Value = "  -6.62%  "
$ws.Cells.Item(36, 4).Value = "0.995"
$ws.Cells.Item(36, 5).Value = "  -0.20%  "
$ws.Cells.Item(37, 4).Value = "17.69"
$ws.Cells.Item(37, 5).Value = "  -2.00%  "
$ws.Cells.Item(38, 5).Value = "  -0.42%  "
$ws.Cells.Item(39, 5).Value = "  -5.82%  "
$ws.Cells.Item(40, 4).Value = "35.83"
$ws.Cells.Item(40, 5).Value = "  -1.50%  "
$ws.Cells.Item(41, 5).Value = "  -6.44%  "
$ws.Cells.Item(42, 4).Value = "0.768"
$ws.Cells.Item(42, 5).Value = "  -6.98%  "
$ws.Cells.Item(43, 4).Value = "128.97"
$ws.Cells.Item(43, 5).Value = "  -1.73%  "
$ws.Cells.Item(44, 5).Value = "  -3.85%  "
$ws.Cells.Item(45, 4).Value = "4.76"
$ws.Cells.Item(45, 5).Value = "  -4.82%  "
$ws.Cells.Item(46, 4).Value = "0.569"
$ws.Cells.Item(46, 5).Value = "  -3.88%  "
$ws.Cells.Item(47, 4).Value = "'252.40"
$ws.Cells.Item(47, 5).Value = "  -7.84%  "
$ws.Cells.Item(48, 4).Value = "0.0895"
$ws.Cells.Item(48, 5).Value = "  -4.12%  "
$ws.Cells.Item(49, 4).Value = "0.0483"
$ws.Cells.Item(49, 5).Value = "  -5.20%  "
$ws.Cells.Item(50, 4).Value = "16.74"
$ws.Cells.Item(50, 5).Value = "  -4.59%  "
$ws.Cells.Item(51, 5).Value = "  -5.20%  "
